$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 489.25
$ws.Range("I18").Value = 489.25
$ws.Range("K18").Value = 489.25
$ws.Range("M18").Value = -205.25

$ws.Range("H32").Value = 1995.5
$ws.Range("I32").Value = 1995.5
$ws.Range("K32").Value = 1995.5
$ws.Range("M32").Value = -1669.5

$ws.Range("H40").Value = 2338.125
$ws.Range("J40").Value = 3334
$ws.Range("L40").Value = 3334
$ws.Range("N40").Value = -3684

$ws.Range("H51").Value = 7263.8
$ws.Range("I51").Value = 4106.3335
$ws.Range("K51").Value = 4106.3335
$ws.Range("M51").Value = -3622.3335

$ws.Range("H58").Value = 5116
$ws.Range("I58").Value = 6032
$ws.Range("J58").Value = 4200
$ws.Range("K58").Value = 18096
$ws.Range("L58").Value = 12600
$ws.Range("M58").Value = -17946
$ws.Range("N58").Value = -12900

$ws.Range("H74").Value = 7180.7144
$ws.Range("I74").Value = 5140
$ws.Range("K74").Value = 5140
$ws.Range("M74").Value = -4204

$ws.Range("H76").Value = 3489.6
$ws.Range("J76").Value = 3499.5
$ws.Range("L76").Value = 3499.5
$ws.Range("N76").Value = -4129.5

$ws.Range("H77").Value = 7180.7144
$ws.Range("I77").Value = 5140
$ws.Range("K77").Value = 25700
$ws.Range("M77").Value = -21020

$ws.Range("H79").Value = 3489.6
$ws.Range("J79").Value = 3499.5
$ws.Range("L79").Value = 3499.5
$ws.Range("N79").Value = -5683.5

$ws.Range("H86").Value = 2763
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 1526
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 1526
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -3772

$ws.Range("H88").Value = 655.5625
$ws.Range("I88").Value = 439
$ws.Range("K88").Value = 439
$ws.Range("M88").Value = -33

$ws.Range("H89").Value = 2763
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 1526
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 7630
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -18862

$ws.Range("H91").Value = 655.5625
$ws.Range("I91").Value = 439
$ws.Range("K91").Value = 439
$ws.Range("M91").Value = 965

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("N93").Value = 0
$ws.Range("L93").ClearContents()

$ws.Range("H98").Value = 3830.3
$ws.Range("I98").Value = 3830.3
$ws.Range("K98").Value = 3830.3
$ws.Range("M98").Value = -2332.3

$ws.Range("H112").Value = 2869.077
$ws.Range("J112").Value = 3419.8
$ws.Range("L112").Value = 10259.4
$ws.Range("N112").Value = -12475.4

$ws.Range("H113").Value = 4799.4
$ws.Range("I113").Value = 4798.5
$ws.Range("K113").Value = 4798.5
$ws.Range("M113").Value = -1544.5

$ws.Range("H117").Value = 40000
$ws.Range("J117").Value = 40000
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

$ws.Range("H122").Value = 3830.3
$ws.Range("I122").Value = 3830.3
$ws.Range("K122").Value = 11490.9
$ws.Range("M122").Value = -9040.900000000001

$ws.Range("H137").Value = 5108051
$ws.Range("I137").Value = 10419579
$ws.Range("K137").Value = 31258737
$ws.Range("M137").Value = -31256187

$ws.Range("H138").Value = 6500
$ws.Range("J138").Value = 7500
$ws.Range("L138").Value = 22500
$ws.Range("N138").Value = -32780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5821
$ws.Range("I61").Value = 3568.5356
$ws.Range("K61").Value = 3568.5356
$ws.Range("M61").Value = -3356.5356

$ws.Range("H122").Value = 3863.5
$ws.Range("I122").Value = 4196.3687
$ws.Range("J122").Value = 2598.6
$ws.Range("K122").Value = 12589.1061
$ws.Range("L122").Value = 7795.799999999999
$ws.Range("M122").Value = -10139.1061
$ws.Range("N122").Value = -12695.8

$ws.Range("H136").Value = 5821
$ws.Range("I136").Value = 3568.5356
$ws.Range("K136").Value = 10705.6068
$ws.Range("M136").Value = -8155.606800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 83336440
$ws.Range("I105").Value = 111113864
$ws.Range("K105").Value = 111113864
$ws.Range("M105").Value = -111112117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 30000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H61").Value = 30000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H138").Value = 120000
$ws.Range("J138").Value = 120000
$ws.Range("L138").Value = 120000
$ws.Range("N138").Value = -130280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6186084
$ws.Range("J4").Value = 100927.5
$ws.Range("L4").Value = 302782.5
$ws.Range("N4").Value = -303006.5

$ws.Range("H57").Value = 692
$ws.Range("I57").Value = 692
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 2076
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = -1517
$ws.Range("M57").ClearContents()

$ws.Range("H132").Value = 1509.5555
$ws.Range("J132").Value = 3002.5
$ws.Range("L132").Value = 27022.5
$ws.Range("N132").Value = -32082.5

$ws.Range("H140").Value = 1162.125
$ws.Range("I140").Value = 959.9
$ws.Range("J140").Value = 1499.1666
$ws.Range("K140").Value = 2879.7
$ws.Range("L140").Value = 4497.4998
$ws.Range("M140").Value = 2300.3
$ws.Range("N140").Value = -14857.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5993.2617
$ws.Range("I132").Value = 4345.4546
$ws.Range("J132").Value = 7805.85
$ws.Range("K132").Value = 13036.3638
$ws.Range("L132").Value = 23417.55
$ws.Range("M132").Value = -10506.3638
$ws.Range("N132").Value = -28477.55

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26284.646
$ws.Range("I7").Value = 29204.777
$ws.Range("J7").Value = 22999.5
$ws.Range("K7").Value = 29204.777
$ws.Range("L7").Value = 22999.5
$ws.Range("M7").Value = -29092.777
$ws.Range("N7").Value = -23223.5

$ws.Range("H22").Value = 2703.0637
$ws.Range("I22").Value = 1750.5862
$ws.Range("J22").Value = 4237.6113
$ws.Range("K22").Value = 1750.5862
$ws.Range("L22").Value = 4237.6113
$ws.Range("M22").Value = -1455.5862
$ws.Range("N22").Value = -4827.6113

$ws.Range("H27").Value = 2703.0637
$ws.Range("I27").Value = 1750.5862
$ws.Range("J27").Value = 4237.6113
$ws.Range("K27").Value = 1750.5862
$ws.Range("L27").Value = 4237.6113
$ws.Range("M27").Value = -1643.5862
$ws.Range("N27").Value = -4451.6113

$ws.Range("H40").Value = 10864.143
$ws.Range("I40").Value = 10554.637
$ws.Range("K40").Value = 10554.637
$ws.Range("M40").Value = -10418.637

$ws.Range("H46").Value = 7256.1875
$ws.Range("I46").Value = 1750
$ws.Range("J46").Value = 9091.583
$ws.Range("K46").Value = 1750
$ws.Range("L46").Value = 9091.583
$ws.Range("M46").Value = -1562
$ws.Range("N46").Value = -9467.583

$ws.Range("H54").Value = 35000
$ws.Range("J54").Value = 35000
$ws.Range("L54").Value = 35000
$ws.Range("N54").Value = -36288

$ws.Range("H55").Value = 1119.4117
$ws.Range("I55").Value = 909.3571
$ws.Range("K55").Value = 909.3571
$ws.Range("M55").Value = -736.3571

$ws.Range("H93").Value = 3161.0557
$ws.Range("I93").Value = 8790
$ws.Range("J93").Value = 996.0769
$ws.Range("K93").Value = 8790
$ws.Range("L93").Value = 996.0769
$ws.Range("M93").Value = -7542
$ws.Range("N93").Value = -3492.0769

$ws.Range("H100").Value = 8336926
$ws.Range("I100").Value = 13160831
$ws.Range("J100").Value = 4726.636
$ws.Range("K100").Value = 13160831
$ws.Range("L100").Value = 4726.636
$ws.Range("M100").Value = -13160290
$ws.Range("N100").Value = -5808.636

$ws.Range("H122").Value = 4979.8
$ws.Range("I122").Value = 5224.75
$ws.Range("K122").Value = 15674.25
$ws.Range("M122").Value = -13224.25

$ws.Range("H126").Value = 26284.646
$ws.Range("I126").Value = 29204.777
$ws.Range("J126").Value = 22999.5
$ws.Range("K126").Value = 87614.33099999999
$ws.Range("L126").Value = 68998.5
$ws.Range("M126").Value = -85144.33099999999
$ws.Range("N126").Value = -73938.5

$ws.Range("H132").Value = 4087.7046
$ws.Range("I132").Value = 3698.3809
$ws.Range("J132").Value = 4443.174
$ws.Range("K132").Value = 11095.1427
$ws.Range("L132").Value = 13329.522
$ws.Range("M132").Value = -8565.1427
$ws.Range("N132").Value = -18389.522

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 27798.8
$ws.Range("J94").Value = 25248.5
$ws.Range("L94").Value = 25248.5
$ws.Range("N94").Value = -27050.5

$ws.Range("H107").Value = 3059.65
$ws.Range("I107").Value = 3321.1765
$ws.Range("J107").Value = 1577.6666
$ws.Range("K107").Value = 9963.5295
$ws.Range("L107").Value = 4732.9998
$ws.Range("M107").Value = -8043.529500000001
$ws.Range("N107").Value = -8572.9998
